$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy column N (rows 3-14, which carries the per-row styling for this table)
# into the new column O, preserving values + formatting, then overwrite the
# values that differ for the new 2021 data column.
$ws.Range("N3:N14").Copy($ws.Range("O3:O14"))

# Row 4: year header
$ws.Range("O4").Value = 2021

# Data rows 5-13
$ws.Range("O5").Value = 97
$ws.Range("O6").Value = 96.2
$ws.Range("O7").Value = 62.7
$ws.Range("O8").Value = 100
$ws.Range("O9").Value = 100
# O10 keeps the copied "-" (no data) marker from N10
$ws.Range("O11").Value = 100
$ws.Range("O12").Value = 57.9
$ws.Range("O13").Value = 100
# O14 keeps the copied "-" (no data) marker from N14

# Update the selection to match the post-edit cursor position
[void]$ws.Range("O17").Select()

Write-Host "done"
